{"js": "// Apply updated salvage/loss table values.\n// Table layout (0-indexed rows incl. header, 0-indexed cols):\n//   row 0 = header\n//   row 1 = Mar 08 ... row 7 = Mar 14\n//   col 0 Date, col 1 Steelhead Daily Salvage, col 2 Steelhead 7-day rolling sum salvage,\n//   col 3 Steelhead Daily Trigger, col 4 Winter-run Daily Loss,\n//   col 5 Winter-run 7-day rolling sum loss, col 6 Winter-run Daily Threshold,\n//   col 7 Winter-run Daily Trigger\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// [rowIndex, colIndex, newText]\nconst edits = [\n  [1, 4, \"2.6\"],   // Mar 08, Winter-run Daily Loss: 2.6* -> 2.6\n  [3, 1, \"4\"],     // Mar 10, Steelhead Daily Salvage: 0 -> 4\n  [3, 2, \"16\"],    // Mar 10, Steelhead 7-day rolling sum salvage: 12 -> 16\n  [4, 1, \"7\"],     // Mar 11, Steelhead Daily Salvage: 0 -> 7\n  [4, 2, \"23\"],    // Mar 11, Steelhead 7-day rolling sum salvage: 12 -> 23\n  [5, 2, \"23\"],    // Mar 12, Steelhead 7-day rolling sum salvage: 12 -> 23\n  [6, 2, \"19\"],    // Mar 13, Steelhead 7-day rolling sum salvage: 8 -> 19\n  [6, 4, \"5.2*\"],  // Mar 13, Winter-run Daily Loss: 0 -> 5.2*\n  [6, 5, \"7.8\"],   // Mar 13, Winter-run 7-day rolling sum loss: 2.6 -> 7.8\n  [7, 2, \"19\"],    // Mar 14, Steelhead 7-day rolling sum salvage: 8 -> 19\n  [7, 5, \"7.8\"],   // Mar 14, Winter-run 7-day rolling sum loss: 2.6 -> 7.8\n];\n\nfor (const [rowIndex, colIndex, newText] of edits) {\n  const cell = table.getCell(rowIndex, colIndex);\n  const para = cell.body.paragraphs.getFirst();\n  const range = para.getRange();\n  range.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Apply updated salvage/loss table values.\n# Table layout (1-indexed rows incl. header, 1-indexed cols, COM style):\n#   row 1 = header\n#   row 2 = Mar 08 ... row 8 = Mar 14\n#   col 1 Date, col 2 Steelhead Daily Salvage, col 3 Steelhead 7-day rolling sum salvage,\n#   col 4 Steelhead Daily Trigger, col 5 Winter-run Daily Loss,\n#   col 6 Winter-run 7-day rolling sum loss, col 7 Winter-run Daily Threshold,\n#   col 8 Winter-run Daily Trigger\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$edits = @(\n    @(2, 5, \"2.6\"),   # Mar 08, Winter-run Daily Loss: 2.6* -> 2.6\n    @(4, 2, \"4\"),     # Mar 10, Steelhead Daily Salvage: 0 -> 4\n    @(4, 3, \"16\"),    # Mar 10, Steelhead 7-day rolling sum salvage: 12 -> 16\n    @(5, 2, \"7\"),     # Mar 11, Steelhead Daily Salvage: 0 -> 7\n    @(5, 3, \"23\"),    # Mar 11, Steelhead 7-day rolling sum salvage: 12 -> 23\n    @(6, 3, \"23\"),    # Mar 12, Steelhead 7-day rolling sum salvage: 12 -> 23\n    @(7, 3, \"19\"),    # Mar 13, Steelhead 7-day rolling sum salvage: 8 -> 19\n    @(7, 5, \"5.2*\"),  # Mar 13, Winter-run Daily Loss: 0 -> 5.2*\n    @(7, 6, \"7.8\"),   # Mar 13, Winter-run 7-day rolling sum loss: 2.6 -> 7.8\n    @(8, 3, \"19\"),    # Mar 14, Steelhead 7-day rolling sum salvage: 8 -> 19\n    @(8, 6, \"7.8\")    # Mar 14, Winter-run 7-day rolling sum loss: 2.6 -> 7.8\n)\n\nforeach ($edit in $edits) {\n    $rowIndex = $edit[0]\n    $colIndex = $edit[1]\n    $newText = $edit[2]\n    $cell = $t.Cell($rowIndex, $colIndex)\n    $cell.Range.Text = $newText\n}\n"}
